# Inserts a new weekly data row at row 51 (Provincia de Cautín, 10/07/2023)
# in the "Rabanito" consolidated sheet, pushing all subsequent rows down by
# one (old row 51 -> new row 52 ... old row 119 -> new row 120).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Shift rows 51..119 down by one row, inserting a fresh blank row at 51.
$ws.Range("A51:R51").Insert("Down")

# Populate the newly inserted row 51 with the new record.
$ws.Range("A51").Value = 10
$ws.Range("B51").Value = "Vega Modelo de Temuco"
$ws.Range("C51").Value = "La Araucanía"
$ws.Range("D51").Value = 45117
$ws.Range("E51").Value = 9
$ws.Range("F51").Value = 300000001
$ws.Range("G51").Value = "Rabanito"
$ws.Range("H51").Value = "Sin especificar"
$ws.Range("I51").Value = "Primera"
$ws.Range("J51").Value = 100
$ws.Range("K51").Value = 8000
$ws.Range("L51").Value = 8000
$ws.Range("M51").Value = 8000
$ws.Range("N51").Value = '$/docena de paquetes'
$ws.Range("O51").Value = "Provincia de Cautín"
$ws.Range("P51").Value = 667
$ws.Range("Q51").Value = 12
$ws.Range("R51").Value = "Hortaliza"
